$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.668.42"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "1.897.37"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5258"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.90%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3806"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07236"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9023"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07634"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.878.25"
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.430"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.66"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008669"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.0000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").Value = "27.701.08"
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").Value = "2.113.71"
$ws.Range("E22").Value = "  +1.09%  "
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.598"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.863"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.187"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.836"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.807"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09140"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05266"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.118"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.70%  "
$ws.Range("E35").Value = "  -1.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7715"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02084"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.562"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.076"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5578"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.74%  "
$ws.Range("E41").Value = "  -0.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.728"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "116.95"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.686"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1509"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4807"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9995"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.598"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "66.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "37.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.22%  "
